$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated / newly-added numeric values per row, matching the target diff.
# Each row array is a list of (ColumnLetter, Value) pairs.

$data = @{
    2  = @{ B = -0.5616080510579985;  C = 0.1593188804880037;  D = -0.549633044125851;
             F = 0.4122776954696132;  G = 0.7251028904350592;  H = -0.1532996908165208;
             I = -0.8832117059949898; J = 0.7789673903946376;  K = 0.1597481019993938 }
    3  = @{ B = 0.1417647591280393;   C = -0.5671871654858154;
             E = 0.3947235741096488;  F = 0.7075487690750948;  G = -0.1708538121764852;
             H = -0.9007658273549541; I = 0.7614132690346732;  J = 0.1421939806394294;
             K = 0.379497744259143 }
    4  = @{ B = -0.4790798465348092;
             D = 0.482830893060655;   E = 0.795656088026101;   F = -0.082746493225479;
             G = -0.8126585084039479; H = 0.8495205879856794;  I = 0.2303012995904356;
             J = 0.4676050632101492 }
    5  = @{ C = 0.6614055265484386;   D = 0.9742307215138846;  E = 0.0958281402623046;
             F = -0.6340838749161644; G = 1.028095221473463;   H = 0.4088759330782192;
             I = 0.6461796966979327;  K = -0.2486961005069136 }
    6  = @{ B = 1.573432754301089;    C = 1.886257949266535;   D = 1.007855368014955;
             E = 0.2779433528364856;  F = 1.940122449226113;   G = 1.320903160830869;
             H = 1.558206924450583;   J = 0.6633311272457364;  K = 1.273820034913197 }
    7  = @{ B = 0.9422837133007778;   C = 0.06388113204919779; D = -0.6660308831292712;
             E = 0.9961482132603562;  F = 0.3769289248651124;  G = 0.6142326884848259;
             I = -0.2806431087200204; J = 0.3298457989474406;  K = 0.1683237681281231 }
    8  = @{ B = 0.0678490295623069;   C = -0.6620629856161621; D = 1.000116110773465;
             E = 0.3808968223782215;  F = 0.6182005859979351;  H = -0.2766752112069113;
             I = 0.3338136964605497;  J = 0.1722916656412322 }
    9  = @{ B = -0.5264228954459207;  C = 1.135756200943707;   D = 0.5165369125484629;
             E = 0.7538406761681764;  G = -0.1410351210366699; H = 0.4694537866307911;
             I = 0.3079317558114735 }
    10 = @{ B = 0.8949500190880419;   C = 0.2757307306927982;  D = 0.5130344943125118;
             F = -0.3818413028923346; G = 0.2286476047751264;  H = 0.06712557395580883 }
    11 = @{ B = 0.2303995154407018;   C = 0.4677032790604154;
             E = -0.427172518144431;  F = 0.18331638952303;    G = 0.02179435870371246 }
    12 = @{ B = 0.4008418571243615;
             D = -0.4940339400804848; E = 0.1164549675869761;  F = -0.04506706323234141 }
    13 = @{ C = -0.5236201424372015;  D = 0.08686876523025952; E = -0.07465326558905801 }
    14 = @{ B = -0.5417707991668423;  C = 0.06871810850061863; D = -0.0928039223186989 }
    15 = @{ B = 0.0506862842519193;   C = -0.1108357465673982 }
    16 = @{ B = -0.1624199859130616 }
}

foreach ($rowKey in $data.Keys) {
    $rowData = $data[$rowKey]
    foreach ($colKey in $rowData.Keys) {
        $cellAddr = "$colKey$rowKey"
        $ws.Range($cellAddr).Value = $rowData[$colKey]
    }
}

$wb.Save()
